$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 8, shifting existing rows 8-12 down to 9-13
$ws.Rows("8").Insert()

# Fill in the newly inserted row 8 with the new data record
$ws.Range("A8").Value = 12
$ws.Range("B8").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C8").Value = "Metropolitana"
$ws.Range("D8").Value = 44467
$ws.Range("E8").Value = 13
$ws.Range("F8").Value = 100112013
$ws.Range("G8").Value = "Alcachofa"
$ws.Range("H8").Value = "Española"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 35
$ws.Range("K8").Value = 12000
$ws.Range("L8").Value = 12000
$ws.Range("M8").Value = 12000
$ws.Range("N8").Value = "$/caja 30 unidades"
$ws.Range("O8").Value = "Provincia de Limarí"
$ws.Range("P8").Value = 400
$ws.Range("Q8").Value = 30
$ws.Range("R8").Value = "Hortaliza"
